# Edit the DBmodel workbook: restructure the "media" table (column B) to add
# title/subtitle/subtype fields (replacing the old name/index/episodes rows at
# that position), drop the "seasons" field from the "collections" table
# (column C), and shift the rest of the media fields down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column B ("media" table) ---
$ws.Cells.Item(3, 2).Value = "title"
$ws.Cells.Item(4, 2).Value = "subtitle"
$ws.Cells.Item(5, 2).Value = "type"
$ws.Cells.Item(6, 2).Value = "subtype"
$ws.Cells.Item(7, 2).Value = "index"
$ws.Cells.Item(8, 2).Value = "episodes"
$ws.Cells.Item(9, 2).Value = "genre"
$ws.Cells.Item(10, 2).Value = "author"
$ws.Cells.Item(11, 2).Value = "status"
$ws.Cells.Item(12, 2).Value = "rating"
$ws.Cells.Item(13, 2).Value = "releaseDate"
$ws.Cells.Item(14, 2).Value = "startDate"
$ws.Cells.Item(15, 2).Value = "completeDate"
$ws.Cells.Item(16, 2).Value = "lastDate"
$ws.Cells.Item(17, 2).Value = "collection_id"

# --- column C ("collections" table): remove the "seasons" field ---
$ws.Cells.Item(5, 3).Clear()

# Match the cell selection left behind in the saved file.
$ws.Range("D4").Select() | Out-Null
